# Rename the monthly journal-entry sheets from abbreviated numeric names
# (e.g. "4_2022") to full month-name form (e.g. "APRIL_2022").

$wb = $excel.ActiveWorkbook

$renameMap = @{
    "4_2022"  = "APRIL_2022"
    "5_2022"  = "MAY_2022"
    "6_2022"  = "JUNE_2022"
    "7_2022"  = "JULY_2022"
    "9_2022"  = "SEPTEMBER_2022"
    "10_2022" = "OCTOBER_2022"
    "11_2022" = "NOVEMBER_2022"
    "12_2022" = "DECEMBER_2022"
    "1_2023"  = "JANUARY_2023"
    "2_2023"  = "FEBRUARY_2023"
    "3_2023"  = "MARCH_2023"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renameMap.ContainsKey($oldName)) {
        $ws.Name = $renameMap[$oldName]
    }
}
